$d = $word.ActiveDocument

function Replace-Range([string]$old, [string]$new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $old"
    }
    $rng.Text = $new
}

# Title
Replace-Range "Cosmic Revelations: A Journey Through Stars and Galaxies" "Exploring the Universe Within: The Fascinating World of Chemistry"

# Author name
Replace-Range " Alana Harper" " Elizabeth Carter"

# Email local part / domain
Replace-Range "alana" "elizabeth"
Replace-Range "harper@astrosciences" "carter@educationalhub"

# Intro paragraph
Replace-Range "The vast expanse of the universe has captivated humanity for ages, stoking our curiosity and igniting our imaginations" "Chemistry, the study of the composition and behavior of matter, offers a thrilling journey into the microscopic realm, where atoms interact to form molecules, shaping the world as we know it"

Replace-Range " As we embark on a journey through the celestial realm, we encounter cosmic wonders that challenge our understanding of space and time. From the fiery hearts of stars to the awe-inspiring grandeur of galaxies, we witness the intricate workings of the universe and unravel the mysteries that lie beyond our earthly realm" " In this essay, we'll unravel the enigma of chemistry, unveiling the secrets of substances, their properties, and the intricate dance of reactions that orchestrate the symphony of life itself"

# Section about stars -> elements
Replace-Range "Beyond our planet's atmosphere, billions of stars, each a celestial beacon of energy and light, twinkle in the night sky" "When we delve into the periodic table, we encounter a cast of elements, each with a unique personality and story to tell"

Replace-Range " We marvel at the fiery dance of these cosmic furnaces, powered by nuclear reactions that sculpt their existence and determine their destiny" " We'll discover how elements combine to form compounds, creating a tapestry of diversity that defines the essence of matter"

Replace-Range " From red dwarfs, the smallest and coolest, to blue supergiants, the largest and most luminous, stars exhibit a mesmerizing array of sizes, colors, and life cycles" " From the mundane to the extraordinary, chemistry is found in every corner of our existence"

Replace-Range " Understanding the intricacies of stellar evolution, we gain insights into the formation and fate of our own Sun and the future of our solar system" " Our bodies, the food we eat, the clothes we wear, the medicines that heal us, and the materials that shape our technological marvels--all bear the imprint of chemistry"

# Section about galaxies -> chemical reactions
Replace-Range "As we traverse the celestial tapestry, we encounter galaxies, vast cosmic cities housing countless stars, gas, and dust" "The study of chemical reactions captivates our imagination, as we witness the transformation of substances into new entities, releasing energy or absorbing it, like a cosmic dance that sculpts the universe"

Replace-Range " From majestic spirals, with their graceful arms swirling like celestial dancers, to elliptical galaxies, harboring ancient secrets within their smooth, elliptical shapes, these celestial metropolises showcase the diversity and grandeur of the universe" " We'll explore the concepts of chemical bonding, energy changes, and equilibrium, gaining insight into the driving forces behind these transformations"

Replace-Range " By studying galaxies, we explore the forces that shape their structure and evolution, unravel the mysteries of dark matter and energy, and glimpse into the vastness of space and time" " From the explosive combustion of fireworks to the subtle interplay of enzymes in our cells, chemistry weaves a symphony of change, a symphony of creation and transformation"

# Summary section
Replace-Range "The journey through the stars and galaxies reveals the boundless wonders of the universe, igniting our curiosity and expanding our understanding of space and time" "In this essay, we embarked on a journey into the realm of chemistry, unraveling the mysteries of matter, elements, compounds, and reactions"

Replace-Range " From the fiery hearts of stars to the vast expanse of galaxies, we uncover the intricate workings of the cosmos and the mysteries that lie beyond our earthly realm" " We explored the periodic table, delved into the concepts of bonding, energy changes, and equilibrium, and witnessed the symphony of transformations that chemistry orchestrates"

Replace-Range " The study of stellar evolution, galaxies, and the forces that govern them grants us insights into the formation and destiny of our solar system, the nature of dark matter and energy, and the vastness of the universe. This exploration of cosmic revelations challenges our preconceptions and invites us to ponder the boundless wonders that lie within the fabric of our universe" " Through this exploration, we gained a deeper understanding of the world around us, and the intricate dance of molecules that shapes our existence"

# Add a trailing empty paragraph at the end of the document body
$d.Content.InsertParagraphAfter()
